# Updates the "Estado de Cuenta" sheet:
#  - Removes the two rows for LIZETH PAOLA HERRERA MUNIZ (periods 2001 and 1912)
#  - Keeps ANA MERCEDES PEREZ DE MARTINEZ's two rows, reordered (2306 then 2307),
#    with the final data row taking on the table's closing bottom-border style
#  - Updates the summary totals (Valor Mora, Cant. Trabajadores, Cant. Periodos)
#  - The signature block rows shift up automatically as rows are deleted

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Grab the closing "last row" border formatting (currently on row 19, LIZETH's
#    final period) and stamp it onto row 17 (which will become the new last data
#    row for ANA once the LIZETH rows are removed).
$ws.Range("B19:J19").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# 2) Swap the period/amount values between row 16 and row 17 so ANA's rows read
#    2306 first, then 2307 (matching the new data order).
$tmpPeriod = $ws.Range("E16").Value2
$tmpMora = $ws.Range("F16").Value2
$tmpSalario = $ws.Range("G16").Value2

$ws.Range("E16").Value2 = $ws.Range("E17").Value2
$ws.Range("F16").Value2 = $ws.Range("F17").Value2
$ws.Range("G16").Value2 = $ws.Range("G17").Value2

$ws.Range("E17").Value2 = $tmpPeriod
$ws.Range("F17").Value2 = $tmpMora
$ws.Range("G17").Value2 = $tmpSalario

# 3) Remove LIZETH's two rows entirely (18 and 19); this shifts the signature
#    block (previously rows 24-25) up to rows 22-23.
$ws.Range("A18:J19").EntireRow.Delete()

# 4) Refresh the summary figures at the top of the statement.
$ws.Range("E11").Value2 = 86614
$ws.Range("C13").Value2 = 1
$ws.Range("F13").Value2 = 2
